$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AppControl")

# Update Run Flag (column B) values on the suite summary rows.
$ws.Range("B3").Value = "Y"
$ws.Range("B7").Value = "N"
$ws.Range("B8").Value = "Y"
$ws.Range("B9").Value = "Y"
$ws.Range("B10").Value = "Y"
$ws.Range("B11").Value = "Y"

# B11 picks up the same formatting as the rest of the column (font size 12).
$ws.Range("B11").Font.Size = 12

# Update the active selection to match the latest edit location.
$ws.Range("B8:B11").Select()
